# Update cryptocurrency price values in column D (Price) to reflect the
# latest scraped data, as produced by the "Updated symbol list" GitHub
# Actions workflow run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Map of row number -> new Price (column D) text value.
# Values are written as text (matching the workbook's existing storage of
# these figures as strings) so that trailing zeros / significant digits
# are preserved exactly instead of being renormalized as numbers.
$updates = [ordered]@{
    2  = "272.09"
    3  = "23.11"
    4  = "6.374"
    5  = "0.06299"
    6  = "3.659"
    7  = "6.758"
    8  = "1.390"
    9  = "0.8337"
    10 = "0.1623"
    11 = "0.08393"
    12 = "0.03480"
    13 = "0.03103"
    14 = "0.09313"
    15 = "3.952"
    16 = "0.001722"
    17 = "0.04881"
    18 = "0.006230"
    23 = "2.321"
    24 = "0.01389"
    27 = "0.0002682"
    40 = "0.04691"
    41 = "0.006904"
    42 = "0.1177"
    43 = "0.003348"
    44 = "0.01252"
    45 = "0.00006269"
    46 = "0.00000000750"
    47 = "0.7889"
    48 = "0.1136"
    49 = "0.00002100"
    50 = "0.01240"
}

foreach ($row in $updates.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    # Prefix with an apostrophe so Excel stores the value as text rather
    # than reinterpreting/rounding it as a number (important for values
    # like "1.390" or "0.00000000750" that carry significant trailing
    # zeros).
    $cell.Value = "'" + $updates[$row]
    # Restore the default "Normal" cell style so no residual text/quote
    # number-format styling is left behind on the cell.
    $cell.Style = "Normal"
}
